$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ver-Construção1")

$ws.Range("D6").Value = "Sim"
$ws.Range("D8").Value = "Sim"
$ws.Range("D10").Value = "NA"
$ws.Range("D11").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D14").Value = "NA"

$ws1 = $wb.Worksheets.Item("Indicadores")
$cos = $ws1.ChartObjects()
$co = $cos.Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection(1)
$sc.Values = @(1,1,1,0)
Write-Host "Values:" $sc.Values
Write-Host "Formula:" $sc.Formula
